$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.555.57"
$ws.Range("E2").Value = "  +2.04%  "
$ws.Range("D3").Value = "1.688.54"
$ws.Range("E4").Value = "  -0.23%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "221.00"
$ws.Range("E5").Value = "  +2.84%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.523"
$ws.Range("E6").Value = "  +0.37%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("E7").Value = "  -0.19%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "30.76"
$ws.Range("E8").Value = "  +3.89%  "
$ws.Range("E9").Value = "  +2.24%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0625"
$ws.Range("E10").Value = "  +1.84%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0901"
$ws.Range("E11").Value = "  -1.64%  "
$ws.Range("D12").Value = "1.932.30"
$ws.Range("E12").Value = "  +3.55%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.78"
$ws.Range("E13").Value = "  +14.12%  "
$ws.Range("B14").Value = "Polygon"
$ws.Range("C14").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.622"
$ws.Range("E14").Value = "  +8.73%  "
$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").Value = "1.686.13"
$ws.Range("E15").Value = "  +3.31%  "
$ws.Range("E16").Value = "  +3.12%  "
$ws.Range("D17").Value = "30.549.01"
$ws.Range("E17").Value = "  +1.99%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "66.01"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "249.65"
$ws.Range("E19").Value = "  -0.02%  "
$ws.Range("D20").Value = "0.0₃0719"
$ws.Range("E20").Value = "  +1.76%  "
$ws.Range("E21").Value = "  -0.26%  "
$ws.Range("E22").Value = "  +3.21%  "
$ws.Range("E23").Value = "  +5.61%  "
$ws.Range("E24").Value = "  +4.66%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "157.79"
$ws.Range("E25").Value = "  -1.23%  "
$ws.Range("E26").Value = "  +1.50%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.112"
$ws.Range("E27").Value = "  +0.45%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.79"
$ws.Range("E28").Value = "  +2.53%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.998"
$ws.Range("E29").Value = "  -0.28%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0500"
$ws.Range("E30").Value = "  +2.17%  "
$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.51"
$ws.Range("E31").Value = "  +3.83%  "
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.14"
$ws.Range("E32").Value = "  +0.77%  "
$ws.Range("B33").Value = "Maker"
$ws.Range("C33").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D33").Value = "1.512.03"
$ws.Range("E33").Value = "  +5.71%  "
$ws.Range("B34").Value = "InternetComputer(DFINITY)"
$ws.Range("C34").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.30"
$ws.Range("E34").Value = "  +2.81%  "
$ws.Range("E35").Value = "  +5.80%  "
$ws.Range("E36").Value = "  -0.40%  "
$ws.Range("E37").Value = "  +4.33%  "
$ws.Range("E38").Value = "  -4.93%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "79.49"
$ws.Range("E39").Value = "  +8.31%  "
$ws.Range("E40").Value = "  +5.46%  "
$ws.Range("E41").Value = "  +1.35%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.854"
$ws.Range("E42").Value = "  +2.79%  "
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.01"
$ws.Range("E43").Value = "  +1.02%  "
$ws.Range("B44").Value = "Kaspa"
$ws.Range("C44").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0505"
$ws.Range("E44").Value = "  +1.84%  "
$ws.Range("E45").Value = "  -2.44%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.998"
$ws.Range("E46").Value = "  -0.17%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "52.36"
$ws.Range("E47").Value = "  -4.38%  "
$ws.Range("D48").Value = "1.824.03"
$ws.Range("E48").Value = "  +2.86%  "
$ws.Range("E49").Value = "  -0.50%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "95.75"
$ws.Range("E50").Value = "  +6.47%  "
$ws.Range("E51").Value = "  +7.23%  "
